# Updates cryptos list figures (price & 1h volume change) and restores
# two swapped coin-name/link/price rows, per upstream GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Bitcoin)
$ws.Range("D2").Value = "51.107.78"
$ws.Range("E2").Value = "  -0.91%  "

# Row 3 (Ethereum)
$ws.Range("D3").Value = "3.052.28"
$ws.Range("E3").Value = "  +0.87%  "

# Row 4 (TetherUSD)
$ws.Range("E4").Value = "  +0.10%  "

# Row 5 (BNB)
$ws.Range("D5").Value = "'387.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.98%  "

# Row 6 (Solana)
$ws.Range("D6").Value = "'101.60"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.78%  "

# Row 7 (XRP)
$ws.Range("D7").Value = "'0.533"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.08%  "

# Row 8 (USDC)
$ws.Range("E8").Value = "  +0.07%  "

# Row 9 (Cardano)
$ws.Range("D9").Value = "'0.576"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.15%  "

# Row 10 (Avalanche)
$ws.Range("D10").Value = "'36.58"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.16%  "

# Row 11 (TRON)
$ws.Range("E11").Value = "  +0.44%  "

# Row 12 (Dogecoin)
$ws.Range("D12").Value = "'0.0845"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.68%  "

# Row 13 (WrappedliquidstakedEther2.0)
$ws.Range("D13").Value = "3.538.89"
$ws.Range("E13").Value = "  +1.11%  "

# Row 14 (Chainlink)
$ws.Range("D14").Value = "'18.21"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.47%  "

# Row 15 (Polkadot)
$ws.Range("D15").Value = "'7.65"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.97%  "

# Row 16 (WrappedEther)
$ws.Range("D16").Value = "3.052.67"
$ws.Range("E16").Value = "  +0.71%  "

# Row 17 (Polygon)
$ws.Range("D17").Value = "'0.984"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.26%  "

# Row 18 (Uniswap)
$ws.Range("D18").Value = "'10.57"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.49%  "

# Row 19 (WrappedBTC)
$ws.Range("D19").Value = "51.142.26"
$ws.Range("E19").Value = "  -0.83%  "

# Row 20 (ImmutableX)
$ws.Range("D20").Value = "'3.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.79%  "

# Row 21 (InternetComputer(DFINITY))
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").Value = "0.0₃0951"
$ws.Range("E21").Value = "  -0.97%  "

# Row 22 (ShibaInu)
$ws.Range("B22").Value = "InternetComputer(DFINITY)"
$ws.Range("C22").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D22").Value = "'12.17"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.32%  "

# Row 23 (Litecoin)
$ws.Range("D23").Value = "'69.54"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.53%  "

# Row 24 (BitcoinCash)
$ws.Range("D24").Value = "'263.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.57%  "

# Row 25 (PancakeSwap)
$ws.Range("D25").Value = "'3.10"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.34%  "

# Row 26 (Filecoin)
$ws.Range("D26").Value = "'7.84"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -7.49%  "

# Row 27 (EthereumClassic)
$ws.Range("D27").Value = "'26.68"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.02%  "

# Row 28 (RenderToken)
$ws.Range("B28").Value = "Dai"
$ws.Range("C28").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D28").Value = "'0.998"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.24%  "

# Row 29 (Dai)
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").Value = "'7.10"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.48%  "

# Row 30 (Kaspa)
$ws.Range("D30").Value = "'0.162"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.71%  "

# Row 31 (Hedera)
$ws.Range("E31").Value = "  -3.41%  "

# Row 32 (Cosmos)
$ws.Range("D32").Value = "'10.42"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.60%  "

# Row 33 (VeChain)
$ws.Range("D33").Value = "'0.0479"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.05%  "

# Row 34 (InjectiveProtocol)
$ws.Range("D34").Value = "'35.32"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.56%  "

# Row 35 (Toncoin)
$ws.Range("D35").Value = "'2.06"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.25%  "

# Row 36 (OKB)
$ws.Range("D36").Value = "'49.93"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.15%  "

# Row 37 (FirstDigitalUSD)
$ws.Range("E37").Value = "  -0.01%  "

# Row 38 (LidoDAOToken)
$ws.Range("D38").Value = "'3.32"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.29%  "

# Row 39 (TheGraph)
$ws.Range("D39").Value = "'0.290"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.39%  "

# Row 40 (Monero)
$ws.Range("D40").Value = "'129.43"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.29%  "

# Row 41 (Celestia)
$ws.Range("D41").Value = "'16.43"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.24%  "

# Row 42 (ARBITRUM)
$ws.Range("D42").Value = "'1.81"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.10%  "

# Row 43 (Stellar)
$ws.Range("D43").Value = "'0.114"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.44%  "

# Row 44 (NEARProtocol)
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").Value = "'2.47"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.75%  "

# Row 45 (Stacks)
$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").Value = "'3.73"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.39%  "

# Row 46 (EnergySwap)
$ws.Range("D46").Value = "'21.52"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.24%  "

# Row 47 (ApeXProtocol)
$ws.Range("E47").Value = "  +3.81%  "

# Row 48 (WEMIXToken)
$ws.Range("E48").Value = "  -0.56%  "

# Row 49 (Maker)
$ws.Range("D49").Value = "2.056.15"
$ws.Range("E49").Value = "  +1.67%  "

# Row 50 (BEAM)
$ws.Range("D50").Value = "'0.0320"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.82%  "

# Row 51 (Mantle)
$ws.Range("D51").Value = "'0.887"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +11.68%  "
